$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Update Status for fed47d42-... row (row 3) to "Ready for handoff"
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"

$zhcn.Range("B3").Value = "Ready for handoff"
$dede.Range("B3").Value = "Ready for handoff"

# Update Latest Handoff Datetime for 5eec1f23-... rows (row 2 and row 3, same shared value)
$zhcn.Range("D2").Value = "2016-03-01 09:51:07"
$zhcn.Range("D3").Value = "2016-03-01 09:51:07"

$dede.Range("D2").Value = "2016-03-01 09:51:18"
$dede.Range("D3").Value = "2016-03-01 09:51:18"
